$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing weight values (B2:B4) with new save-state from previous weights
$ws.Range("B2").Value = 6.268868452723979
$ws.Range("B3").Value = -0.4033703316514688
$ws.Range("B4").Value = -2.676184643176573

# Remove the now-stale 4th data row (epoch 3 / row 5) entirely, shifting cells up
$ws.Range("A5:C5").Delete()
